$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Formula = "=G2+H2"
$ws.Range("M2").Formula = "=L2*(J2/100)"
$ws.Range("N2").Formula = "=M2*(K2/100)"
